$wb = $excel.ActiveWorkbook

# --- PHENOTYPES sheet: add a helper list of possible data types and
#     turn the "Data Type" column validation into a dropdown list ---
$wsPheno = $wb.Worksheets.Item("PHENOTYPES")

$wsPheno.Range("J1").Value = "Possible data types"
$wsPheno.Range("J2").Value = "float"
$wsPheno.Range("J3").Value = "int"
$wsPheno.Range("J4").Value = "char"

# Replace the old free-text "Data Type" validation with a dropdown list
# driven by the J2:J4 helper range. This keeps old templates compatible
# (importer still just reads text out of column D).
$wsPheno.Range("D1").Validation.Delete()
$wsPheno.Range("D1:D1048576").Validation.Add(3, 1, 1, "=`$J`$2:`$J`$4")

$wsPheno.Activate()
$wsPheno.Range("D2").Select()

# --- LOCATION sheet: move the selection, this is no longer the active tab ---
$wsLoc = $wb.Worksheets.Item("LOCATION")
$wsLoc.Activate()
$wsLoc.Range("C25").Select()

# --- METADATA sheet becomes the active tab again ---
$wsMeta = $wb.Worksheets.Item("METADATA")
$wsMeta.Activate()
$wsMeta.Range("C2").Select()

Write-Output "done"
